$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$msgShort = '⏳ No pierdas más clientes. Contáctanos ahora y asegura tu lugar. ¡Es tu momento de destacar en internet! 🌟 📩 
Escríbenos ya para obtener esta oferta exclusiva. ¡No esperes más! 😊'
$msgLong = '¡Hola! 👋🏼 ¿Estás listo para llevar tu negocio al siguiente nivel? 🌐✨ 🚀 
Transforma tu presencia en línea y atrapa a esos clientes que están buscando lo que ofreces. ¡Imagina un sitio web profesional y atractivo que convierta visitantes en compradores fieles! 🛍️🔥 
¿Listo para comenzar? Responde a este mensaje y uno de nuestros expertos te guiará en el proceso. ¡Es tu momento de brillar en internet! 🌐🌟'
$imgA = 'C:\Users\4to CREATIVO\Desktop\WhatsAutoA\Images\A.jpeg'
$imgB = 'C:\Users\4to CREATIVO\Desktop\WhatsAutoA\Images\B.jpeg'

# --- Update phone numbers (column A) for existing rows 2-26 ---
$ws.Range("A2").Value = "'7712298632"
$ws.Range("A3").Value = "'7711258720"
$ws.Range("A4").Value = "'7717155410"
$ws.Range("A5").Value = "'5583699920"
$ws.Range("A6").Value = "'7715690454"
$ws.Range("A7").Value = "'7711615285"
$ws.Range("A8").Value = "'7711575091"
$ws.Range("A9").Value = "'7713187830"
$ws.Range("A10").Value = "'7717120945"
$ws.Range("A11").Value = "'7711248373"
$ws.Range("A12").Value = "'7711519561"
$ws.Range("A13").Value = "'7711468170"
$ws.Range("A14").Value = "'7711015664"
$ws.Range("A15").Value = "'7711227813"
$ws.Range("A16").Value = "'7717091232"
$ws.Range("A17").Value = "'7713031132"
$ws.Range("A18").Value = "'7713496839"
$ws.Range("A19").Value = "'7711301339"
$ws.Range("A20").Value = "'7711004270"
$ws.Range("A21").Value = "'7717021603"
$ws.Range("A22").Value = "'7711405885"
$ws.Range("A23").Value = "'7715268539"
$ws.Range("A24").Value = "'7712076789"
$ws.Range("A25").Value = "'7712894116"
$ws.Range("A26").Value = "'7711869908"

# --- Append new rows 27-51 ---
$ws.Range("A27").Value = "'7711809278"
$ws.Range("B27").Value = $msgShort
$ws.Range("C27").Value = $imgA

$ws.Range("A28").Value = "'7712141045"
$ws.Range("B28").Value = $msgShort
$ws.Range("C28").Value = $imgA

$ws.Range("A29").Value = "'7712950081"
$ws.Range("B29").Value = $msgLong
$ws.Range("C29").Value = $imgB

$ws.Range("A30").Value = "'7712444491"
$ws.Range("B30").Value = $msgShort
$ws.Range("C30").Value = $imgA

$ws.Range("A31").Value = "'7712994514"
$ws.Range("B31").Value = $msgShort
$ws.Range("C31").Value = $imgA

$ws.Range("A32").Value = "'7712292105"
$ws.Range("B32").Value = $msgLong
$ws.Range("C32").Value = $imgB

$ws.Range("A33").Value = "'7711393462"
$ws.Range("B33").Value = $msgShort
$ws.Range("C33").Value = $imgA

$ws.Range("A34").Value = "'7711372699"
$ws.Range("B34").Value = $msgLong
$ws.Range("C34").Value = $imgB

$ws.Range("A35").Value = "'7717470972"
$ws.Range("B35").Value = $msgLong
$ws.Range("C35").Value = $imgB

$ws.Range("A36").Value = "'7712285031"
$ws.Range("B36").Value = $msgLong
$ws.Range("C36").Value = $imgB

$ws.Range("A37").Value = "'7717724295"
$ws.Range("B37").Value = $msgShort
$ws.Range("C37").Value = $imgA

$ws.Range("A38").Value = "'7712955716"
$ws.Range("B38").Value = $msgLong
$ws.Range("C38").Value = $imgB

$ws.Range("A39").Value = "'7712167964"
$ws.Range("B39").Value = $msgShort
$ws.Range("C39").Value = $imgA

$ws.Range("A40").Value = "'7716992902"
$ws.Range("B40").Value = $msgLong
$ws.Range("C40").Value = $imgB

$ws.Range("A41").Value = "'7711566394"
$ws.Range("B41").Value = $msgShort
$ws.Range("C41").Value = $imgA

$ws.Range("A42").Value = "'7714100900"
$ws.Range("B42").Value = $msgLong
$ws.Range("C42").Value = $imgB

$ws.Range("A43").Value = "'7712052619"
$ws.Range("B43").Value = $msgLong
$ws.Range("C43").Value = $imgB

$ws.Range("A44").Value = "'7712078109"
$ws.Range("B44").Value = $msgShort
$ws.Range("C44").Value = $imgA

$ws.Range("A45").Value = "'7711991969"
$ws.Range("B45").Value = $msgShort
$ws.Range("C45").Value = $imgA

$ws.Range("A46").Value = "'7711393226"
$ws.Range("B46").Value = $msgShort
$ws.Range("C46").Value = $imgA

$ws.Range("A47").Value = "'7721015566"
$ws.Range("B47").Value = $msgShort
$ws.Range("C47").Value = $imgA

$ws.Range("A48").Value = "'7711255859"
$ws.Range("B48").Value = $msgLong
$ws.Range("C48").Value = $imgB

$ws.Range("A49").Value = "'7712180247"
$ws.Range("B49").Value = $msgShort
$ws.Range("C49").Value = $imgA

$ws.Range("A50").Value = "'7712954824"
$ws.Range("B50").Value = $msgShort
$ws.Range("C50").Value = $imgA

$ws.Range("A51").Value = "'7716834318"
$ws.Range("B51").Value = $msgShort
$ws.Range("C51").Value = $imgA

